# Fin du projet, resolution du fetch POST
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (test #2)
$ws.Cells.Item(3, 2).Value = "affichage des produits avec leurs noms et description"
$ws.Cells.Item(3, 4).Value = "Affichage dynamique de chaque produit"
$ws.Cells.Item(3, 5).Value = "OK"

# Row 4 (test #3)
$ws.Cells.Item(4, 2).Value = "Redirection vers le produit choisit"
$ws.Cells.Item(4, 3).Value = "Clic sur un produit"
$ws.Cells.Item(4, 4).Value = "ouvre une nouvelle page avec le bon produit"
$ws.Cells.Item(4, 5).Value = "OK"

# Row 5 (test #4)
$ws.Cells.Item(5, 2).Value = "affichage du choix de couleurs en fonction du produit"
$ws.Cells.Item(5, 3).Value = "clic sur séléctionner une couleur "
$ws.Cells.Item(5, 4).Value = "affiche exactement les nombres de couleurs disponible pour le produit"
$ws.Cells.Item(5, 5).Value = "OK"

# Row 6 (test #5)
$ws.Cells.Item(6, 2).Value = "Ajout au panier"
$ws.Cells.Item(6, 3).Value = "choisir une couleur, une quantité et click ajouter au panier"
$ws.Cells.Item(6, 4).Value = "Ajout du produit dans la page panier"
$ws.Cells.Item(6, 5).Value = "OK"

# Row 7 (test #6) - now has a number instead of the previous "…" text
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Page panier"
$ws.Cells.Item(7, 3).Value = "clic sur panier dans la barre nav"
$ws.Cells.Item(7, 4).Value = "Affiche la page du panier avec les articles séléctionner"
$ws.Cells.Item(7, 5).Value = "OK"

# Row 8 (test #7)
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Changer la quantité d'un produit"
$ws.Cells.Item(8, 3).Value = "choisir une quantité"
$ws.Cells.Item(8, 4).Value = "Augmente ou diminué la quantité et modifie le prix des produits et du panier"
$ws.Cells.Item(8, 5).Value = "OK"

# Row 9 (test #8)
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Supprimer un produit"
$ws.Cells.Item(9, 3).Value = "clic sur supprimer à côté du produit"
$ws.Cells.Item(9, 4).Value = "Suppression du produit"
$ws.Cells.Item(9, 5).Value = "OK"

# Row 10 (test #9)
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Passer commande"
$ws.Cells.Item(10, 3).Value = "remplir le formulaire"
$ws.Cells.Item(10, 4).Value = "validation de la commande et recuperation du numero de commande"
$ws.Cells.Item(10, 5).Value = "OK"

# Row heights: rows 3-8 and row 10 grow to fit the new wrapped text, row 9 stays as-is
$ws.Rows.Item(3).RowHeight = 43.5
$ws.Rows.Item(4).RowHeight = 43.5
$ws.Rows.Item(5).RowHeight = 43.5
$ws.Rows.Item(6).RowHeight = 43.5
$ws.Rows.Item(7).RowHeight = 43.5
$ws.Rows.Item(8).RowHeight = 43.5
$ws.Rows.Item(10).RowHeight = 43.5

# Update the active selection to reflect where the author finished editing
$ws.Range("E10").Select()
